# Add a "False alarm rate" / "false_alarm_rate" column (column N) to each of
# the 7 metrics worksheets, copying the header style from the existing M1
# header cell and filling in the new per-row numeric values.

$wb = $excel.ActiveWorkbook

$headers = @(
    "False alarm rate",
    " false_alarm_rate",
    "false_alarm_rate",
    "false_alarm_rate",
    " false_alarm_rate",
    " false_alarm_rate",
    "false_alarm_rate"
)

$values = @(
    @(0.01446845504222553, 0.01428216592151018, 0.01415797317436662, 0.02608047690014903, 0.6434426229508197),
    @(0.01434426229508197, 0.01428216592151018, 0.01415797317436662, 0.01825633383010432, 0.2174615002483855),
    @(0.02260307998012916, 0.01428216592151018, 0.0142200695479384,  0.03446348733233979, 0.06756085444610035),
    @(0.01403378042722305, 0.01434426229508197, 0.0142200695479384,  0.01484103328365623, 0.04570293094883259),
    @(0.01813214108296075, 0.01353700943864878, 0.01341281669150522, 0.05253353204172876, 0.1058743169398907),
    @(0.0153378042722305,  0.013474913065077,   0.01335072031793343, 0.07153502235469449, 0.06607054148037754),
    @(0.01589667163437655, 0.01341281669150522, 0.01341281669150522, 0.02185792349726776, 0.09227521112767015)
)

for ($i = 1; $i -le 7; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Copy the style of the existing "F1 Score" header (M1) onto the new
    # header cell N1, then overwrite its value/text so it doesn't pick up
    # the copied text of M1.
    $ws.Range("M1").Copy($ws.Range("N1"))
    $ws.Range("N1").Value = $headers[$i - 1]

    $rowVals = $values[$i - 1]
    for ($r = 2; $r -le 6; $r++) {
        $ws.Cells.Item($r, 14).Value = $rowVals[$r - 2]
    }
}
